$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.299544
$ws.Range("H2").Value = 42.898632
$ws.Range("I2").Value = 0.3403136425785068
$ws.Range("J2").Value = 0.3403136425785068
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.056841666666666
$ws.Range("N2").Value = 3.170525
$ws.Range("O2").Value = 0.04271158973519736
$ws.Range("P2").Value = 0.04271158973519736
$ws.Range("Q2").Value = 15.11235391353333
$ws.Range("R2").Value = 136.0111852218
$ws.Range("S2").Value = 0.01453533668310378
$ws.Range("T2").Value = 0.01453533668310378

$ws.Range("G3").Value = 14.299544
$ws.Range("H3").Value = 42.898632
$ws.Range("I3").Value = 0.3403136425785068
$ws.Range("J3").Value = 0.3403136425785068
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 19.80485166666667
$ws.Range("N3").Value = 59.41455500000001
$ws.Range("O3").Value = 0.8004005953144415
$ws.Range("P3").Value = 0.8004005953144415
$ws.Range("Q3").Value = 283.2003478209733
$ws.Range("R3").Value = 2548.80313038876
$ws.Range("S3").Value = 0.2723872421134629
$ws.Range("T3").Value = 0.2723872421134629

$ws.Range("G4").Value = 14.299544
$ws.Range("H4").Value = 42.898632
$ws.Range("I4").Value = 0.3403136425785068
$ws.Range("J4").Value = 0.3403136425785068
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.881981
$ws.Range("N4").Value = 11.645943
$ws.Range("O4").Value = 0.1568878149503611
$ws.Range("P4").Value = 0.1568878149503611
$ws.Range("Q4").Value = 55.51055811666399
$ws.Range("R4").Value = 499.595023049976
$ws.Range("S4").Value = 0.05339106378194011
$ws.Range("T4").Value = 0.0533910637819401

$ws.Range("G5").Value = 19.365057
$ws.Range("H5").Value = 58.095171
$ws.Range("I5").Value = 0.4608673595752713
$ws.Range("J5").Value = 0.4608673595752713
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.056841666666666
$ws.Range("N5").Value = 3.170525
$ws.Range("O5").Value = 0.04271158973519736
$ws.Range("P5").Value = 0.04271158973519736
$ws.Range("Q5").Value = 20.465799114975
$ws.Range("R5").Value = 184.192192034775
$ws.Range("S5").Value = 0.01968437758452267
$ws.Range("T5").Value = 0.01968437758452267

$ws.Range("G6").Value = 19.365057
$ws.Range("H6").Value = 58.095171
$ws.Range("I6").Value = 0.4608673595752713
$ws.Range("J6").Value = 0.4608673595752713
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.80485166666667
$ws.Range("N6").Value = 59.41455500000001
$ws.Range("O6").Value = 0.8004005953144415
$ws.Range("P6").Value = 0.8004005953144415
$ws.Range("Q6").Value = 383.5220814015451
$ws.Range("R6").Value = 3451.698732613906
$ws.Range("S6").Value = 0.3688785089650419
$ws.Range("T6").Value = 0.3688785089650419

$ws.Range("G7").Value = 19.365057
$ws.Range("H7").Value = 58.095171
$ws.Range("I7").Value = 0.4608673595752713
$ws.Range("J7").Value = 0.4608673595752713
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.881981
$ws.Range("N7").Value = 11.645943
$ws.Range("O7").Value = 0.1568878149503611
$ws.Range("P7").Value = 0.1568878149503611
$ws.Range("Q7").Value = 75.17478333791699
$ws.Range("R7").Value = 676.5730500412529
$ws.Range("S7").Value = 0.07230447302570669
$ws.Range("T7").Value = 0.07230447302570668

$ws.Range("G8").Value = 8.35412
$ws.Range("H8").Value = 25.06236
$ws.Range("I8").Value = 0.1988189978462219
$ws.Range("J8").Value = 0.1988189978462219
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.056841666666666
$ws.Range("N8").Value = 3.170525
$ws.Range("O8").Value = 0.04271158973519736
$ws.Range("P8").Value = 0.04271158973519736
$ws.Range("Q8").Value = 8.828982104333331
$ws.Range("R8").Value = 79.46083893899998
$ws.Range("S8").Value = 0.008491875467570919
$ws.Range("T8").Value = 0.008491875467570919

$ws.Range("G9").Value = 8.35412
$ws.Range("H9").Value = 25.06236
$ws.Range("I9").Value = 0.1988189978462219
$ws.Range("J9").Value = 0.1988189978462219
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 19.80485166666667
$ws.Range("N9").Value = 59.41455500000001
$ws.Range("O9").Value = 0.8004005953144415
$ws.Range("P9").Value = 0.8004005953144415
$ws.Range("Q9").Value = 165.4521074055334
$ws.Range("R9").Value = 1489.0689666498
$ws.Range("S9").Value = 0.1591348442359367
$ws.Range("T9").Value = 0.1591348442359367

$ws.Range("G10").Value = 8.35412
$ws.Range("H10").Value = 25.06236
$ws.Range("I10").Value = 0.1988189978462219
$ws.Range("J10").Value = 0.1988189978462219
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.881981
$ws.Range("N10").Value = 11.645943
$ws.Range("O10").Value = 0.1568878149503611
$ws.Range("P10").Value = 0.1568878149503611
$ws.Range("Q10").Value = 32.43053511172
$ws.Range("R10").Value = 291.8748160054799
$ws.Range("S10").Value = 0.03119227814271431
$ws.Range("T10").Value = 0.0311922781427143
